$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.262.82'
$ws.Range("E2").Value = '  -1.39%  '
$ws.Range("D3").Value = '2.248.32'
$ws.Range("E3").Value = '  -1.36%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.66'
$ws.Range("E5").Value = '  -1.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.631'
$ws.Range("E6").Value = '  -1.66%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '74.66'
$ws.Range("E7").Value = '  -5.53%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.620'
$ws.Range("E9").Value = '  -3.93%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.67'
$ws.Range("E10").Value = '  +3.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0948'
$ws.Range("E11").Value = '  -2.41%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.18'
$ws.Range("E12").Value = '  -2.73%  '
$ws.Range("E13").Value = '  -1.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.46'
$ws.Range("E14").Value = '  -4.73%  '
$ws.Range("E15").Value = '  -1.87%  '
$ws.Range("D16").Value = '2.230.67'
$ws.Range("E16").Value = '  -1.85%  '
$ws.Range("D17").Value = '42.150.66'
$ws.Range("E17").Value = '  -1.46%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000101'
$ws.Range("E18").Value = '  +1.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.41'
$ws.Range("E19").Value = '  +0.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.14'
$ws.Range("E20").Value = '  -1.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.23'
$ws.Range("E21").Value = '  +2.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '231.54'
$ws.Range("E22").Value = '  -1.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.95'
$ws.Range("E23").Value = '  +34.49%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.38'
$ws.Range("E25").Value = '  +0.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.61'
$ws.Range("E26").Value = '  -4.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.32'
$ws.Range("E27").Value = '  -1.36%  '
$ws.Range("E28").Value = '  +3.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '171.01'
$ws.Range("E29").Value = '  +1.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.96'
$ws.Range("E30").Value = '  +0.28%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0821'
$ws.Range("E31").Value = '  -3.54%  '
$ws.Range("E32").Value = '  -0.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '30.71'
$ws.Range("E33").Value = '  +0.73%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.126'
$ws.Range("E34").Value = '  -1.68%  '
$ws.Range("E35").Value = '  +9.77%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.46'
$ws.Range("E36").Value = '  -2.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0314'
$ws.Range("E37").Value = '  +3.49%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '13.80'
$ws.Range("E38").Value = '  +1.84%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.18'
$ws.Range("E39").Value = '  -3.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.79'
$ws.Range("E40").Value = '  -2.26%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '62.61'
$ws.Range("E41").Value = '  +1.90%  '
$ws.Range("E42").Value = '  -2.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '107.00'
$ws.Range("E43").Value = '  -7.40%  '
$ws.Range("E44").Value = '  +0.79%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.70'
$ws.Range("E45").Value = '  -2.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.997'
$ws.Range("E46").Value = '  -0.38%  '
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.12'
$ws.Range("E47").Value = '  -2.43%  '
$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.18'
$ws.Range("E48").Value = '  -0.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.30'
$ws.Range("E49").Value = '  +2.10%  '
$ws.Range("E50").Value = '  -10.01%  '
